# "Added ids to Atlanta data, ids to php pull, scores to ATL crimes"
#
# The "Scores" sheet lists crime categories per city (Atlanta in A/B,
# New York City in C/D, Chicago in E/F). This edit:
#   1) Corrects the Atlanta crime-list typo "LARCENY-NONE VEHICLE" ->
#      "LARCENY-NON VEHICLE".
#   2) Fills in the previously-empty "ATL Score" column (B) with the
#      weighting score for each Atlanta crime category.
#   3) Makes "Scores" the active/selected sheet (it had been "City"),
#      with the last selection left on cell E14.

$wb  = $excel.ActiveWorkbook
$ws  = $wb.Worksheets.Item("Scores")

# Make "Scores" the active sheet and leave the cursor on E14.
$ws.Activate() | Out-Null
$ws.Range("E14").Select() | Out-Null

# Fix the Atlanta crime-list label typo.
$ws.Range("A8").Value = "LARCENY-NON VEHICLE"

# Fill in the new "ATL Score" values (column B, rows 2-12) alongside the
# existing Atlanta crime list in column A.
$atlScores = @{
    2  = 1
    3  = 0.25
    4  = 0.25
    5  = 0.5
    6  = 1
    7  = 0.25
    8  = 0.5
    9  = 1
    10 = 0.25
    11 = 1
    12 = 1
}

foreach ($row in $atlScores.Keys) {
    $ws.Cells.Item($row, 2).Value = $atlScores[$row]
}
